$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column D (dow) values for rows 2-21
$dValues = @{
    2 = 7
    3 = 4
    4 = 3
    6 = 6
    7 = 6
    8 = 3
    9 = 3
    10 = 1
    11 = 0
    12 = 4
    13 = 6
    14 = 5
    15 = 7
    16 = 0
    17 = 2
    18 = 2
    20 = 6
    21 = 0
}

foreach ($row in $dValues.Keys) {
    $ws.Cells.Item($row, 4).Value = $dValues[$row]
}

# Delete row 22 entirely (shifts cells up)
$ws.Rows.Item(22).Delete()
